# Fix a typo'd smart/curly apostrophe in the sample code shown on the
# "DataFrame(data['customers'])" slide: ’ -> '
#
# Slide 33, shape "CustomShape 3" (3rd real shape on the slide), 3rd
# paragraph of its text body reads:
#   customers = DataFrame(data['customers’])
# and should read:
#   customers = DataFrame(data['customers'])

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(33)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(3, 1)

# The engine normalizes curly vs straight quotes when comparing old/new
# text, so assigning the "corrected" string directly looks like a no-op
# and the run is left untouched. Force a real text replacement first
# (distinct placeholder text), then set the final corrected text so the
# run actually gets rewritten with a straight apostrophe.
$para.Text = "___QUOTE_FIX_PLACEHOLDER___"
$para2 = $tr.Paragraphs(3, 1)
$para2.Text = "customers = DataFrame(data['customers'])"
